# ---------------------------------------------------------------------------
# Burundi final.xlsx - add VehicleFleet sheet & update Warehouses data
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Warehouses sheet: rename warehouse labels + update Capacity values
# ---------------------------------------------------------------------------
$wsWarehouses = $wb.Worksheets.Item("Warehouses")

$wsWarehouses.Cells.Item(2,1).Value = "BUJUMBURA"
$wsWarehouses.Cells.Item(2,4).Value = 11917

$wsWarehouses.Cells.Item(3,1).Value = "NGOZI"
$wsWarehouses.Cells.Item(3,4).Value = 11917

$wsWarehouses.Cells.Item(4,1).Value = "GITEGA"
$wsWarehouses.Cells.Item(4,4).Value = 11917

$wsWarehouses.Range("A2:A4").NumberFormat = "@"

# Column widths
$wsWarehouses.Columns.Item(4).ColumnWidth = 12.3659

# ---------------------------------------------------------------------------
# 2. Add the VehicleFleet worksheet (last tab)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFleet = $wb.Worksheets.Add($null, $lastSheet)
$wsFleet.Name = "VehicleFleet"

# Header row
$wsFleet.Cells.Item(1,1).Value = "Warehouse"
$wsFleet.Cells.Item(1,2).Value = "Plate Nr"
$wsFleet.Cells.Item(1,3).Value = "Make"
$wsFleet.Cells.Item(1,4).Value = "Model"
$wsFleet.Cells.Item(1,5).Value = "Capacity in MT"

# Data rows: Warehouse, Plate Nr, Make, Model, Capacity in MT
$fleetRows = @(
    @("GITEGA","CD44A95","RENAULT 6X4","350,34",18),
    @("GITEGA","CD44B02","RENAULT 6X4","350,34",18),
    @("GITEGA","CD44A89","RENAULT4X4","300,19",8),
    @("GITEGA","CD44A91","RENAULT 4X4","300,19",8),
    @("GITEGA","CD44A98","RENAULT 4X4","300,19",8),
    @("GITEGA","CD44A54","TOYOTA DYNA",0,3.5),
    @("GITEGA","E059AIT","TOYOTA PIC-UP","Land cruiser",1.5),
    @("GITEGA","CD107-98U","TRAILER",0,15),
    @("BUJUMBURA","CD44A96","RENAULT 6X4","350,34",18),
    @("BUJUMBURA","CD44A52","RENAULT 4X4","300,19",8),
    @("BUJUMBURA","CD44A81","ISUZU",0,4.2),
    @("BUJUMBURA","CD44A55","ISUZU",0,4.2),
    @("BUJUMBURA","CD44A86","ISUZU",0,4.2),
    @("BUJUMBURA","CD44A87","ISUZU",0,4.2),
    @("BUJUMBURA","CD44A35","TOYOTA DYNA",0,3.5),
    @("BUJUMBURA","CD44A25","TOYOTA DYNA",0,3.5),
    @("BUJUMBURA","CD44A31","TOYOTA  PIC-UP","Land cruiser",1.5),
    @("BUJUMBURA","E058AIT","TOYOTA  PIC-UP","Land cruiser",1.5),
    @("NGOZI","CD44A88","RENAULT 6X4","350,34",18),
    @("NGOZI","CD44A94","RENAULT 6X4","350,34",18),
    @("NGOZI","CD44B01","RENAULT 6X4","350,34",18),
    @("NGOZI","CD44A90","RENAULT 4X4","300,19",8),
    @("NGOZI","CD44A57","RENAULT 4X4","300,19",8),
    @("NGOZI","CD44A48","TOYOTA DYNA",0,3.5),
    @("NGOZI","CD44A43","TOYOTA DYNA",0,3.5),
    @("NGOZI","CD44A33","TOYOTA  PIC-UP","Land cruiser",1.5),
    @("NGOZI","CD107-69U","TRAILER",0,15)
)

$r = 2
foreach ($row in $fleetRows) {
    $wsFleet.Cells.Item($r,1).Value = $row[0]
    $wsFleet.Cells.Item($r,2).Value = $row[1]
    $wsFleet.Cells.Item($r,3).Value = $row[2]
    $wsFleet.Cells.Item($r,4).Value = $row[3]
    $wsFleet.Cells.Item($r,5).Value = $row[4]
    $r = $r + 1
}

# Column formats
$wsFleet.Range("A2:A28").NumberFormat = "@"
$wsFleet.Range("E2:E28").NumberFormat = "0.00"

# Column widths
$wsFleet.Columns.Item(1).ColumnWidth = 11.4323
$wsFleet.Columns.Item(3).ColumnWidth = 12.2995
$wsFleet.Columns.Item(5).ColumnWidth = 15.03

# Header alignment
$wsFleet.Range("A1:E1").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 3. Turn the VehicleFleet range into a table
# ---------------------------------------------------------------------------
$fleetTableRange = $wsFleet.Range("A1:E28")
$fleetTable = $wsFleet.ListObjects.Add(1, $fleetTableRange, $null, 1)
$fleetTable.Name = "Tabelle35"
$fleetTable.TableStyle = "TableStyleLight16"

# ---------------------------------------------------------------------------
# 4. Warehouses table restyle
# ---------------------------------------------------------------------------
$whTable = $wsWarehouses.ListObjects.Item(1)
$whTable.TableStyle = "TableStyleLight16"

# ---------------------------------------------------------------------------
# 5. Sheet view / selection updates
# ---------------------------------------------------------------------------
$wsSchools = $wb.Worksheets.Item("Schools")
$wsSchools.Activate()
$wsSchools.Range("F23").Select()

$wsWarehouses.Activate()
$wsWarehouses.Range("A2:A4").Select()

$wsFleet.Activate()
$wsFleet.Range("D29").Select()
